$d = $word.ActiveDocument

$d.Content.Find.Execute("446÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "510÷3=", 2) | Out-Null
$d.Content.Find.Execute("693÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "338÷7=", 2) | Out-Null
$d.Content.Find.Execute("765÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "280÷3=", 2) | Out-Null
$d.Content.Find.Execute("900÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "876÷7=", 2) | Out-Null
$d.Content.Find.Execute("574÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "543÷9=", 2) | Out-Null
$d.Content.Find.Execute("268÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "269÷2=", 2) | Out-Null
$d.Content.Find.Execute("720÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "650÷8=", 2) | Out-Null
$d.Content.Find.Execute("912÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "629÷5=", 2) | Out-Null
$d.Content.Find.Execute("566÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "255÷3=", 2) | Out-Null
$d.Content.Find.Execute("229÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "344÷3=", 2) | Out-Null
$d.Content.Find.Execute("836÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "231÷4=", 2) | Out-Null
$d.Content.Find.Execute("597÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "941÷6=", 2) | Out-Null
$d.Content.Find.Execute("379÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "618÷9=", 2) | Out-Null
$d.Content.Find.Execute("846÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "497÷4=", 2) | Out-Null
$d.Content.Find.Execute("846÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "169÷8=", 2) | Out-Null
$d.Content.Find.Execute("682÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "838÷3=", 2) | Out-Null
$d.Content.Find.Execute("827÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "605÷8=", 2) | Out-Null
$d.Content.Find.Execute("960÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "378÷8=", 2) | Out-Null
$d.Content.Find.Execute("898÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "661÷3=", 2) | Out-Null
$d.Content.Find.Execute("210÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "499÷3=", 2) | Out-Null
$d.Content.Find.Execute("825÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "882÷4=", 2) | Out-Null
$d.Content.Find.Execute("208÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "468÷3=", 2) | Out-Null
$d.Content.Find.Execute("544÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "171÷8=", 2) | Out-Null
$d.Content.Find.Execute("344÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "385÷6=", 2) | Out-Null
$d.Content.Find.Execute("997÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "419÷6=", 2) | Out-Null

Write-Host "Done"
